# Fruta / hortaliza, semanal
# Insert 3 new weekly data rows (Artic Star variety) at the top of the
# Nectarin data block, pushing the existing rows 67-142 down to 70-145.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows before row 67; this shifts existing rows 67-142
# down to 70-145, carrying all of their data/formatting with them and
# updating the sheet dimension automatically.
$ws.Rows("67:69").Insert()

# Common (constant) columns for every data row in this sheet.
$mercadoId = 8
$mercado   = "Terminal La Palmera de La Serena"
$region    = "Coquimbo"
$codreg    = 4
$tipo      = "Fruta"
$productoId = 100103
$producto   = "Frutos de hueso (carozo)"
$categoriaId = 100103006
$categoria   = "Nectarín"

$newRows = @(
    @{ Row=67; Fecha=44539; Variedad="Artic Star"; Calidad="Especial"; Volumen=10;  PMin=485000; PMax=490000; PProm=487500; Unidad="$/bins (420 kilos)"; Origen="Región de O'Higgins"; PKg=1161; KgUnidad=420 },
    @{ Row=68; Fecha=44539; Variedad="Artic Star"; Calidad="Primera";  Volumen=20;  PMin=455000; PMax=460000; PProm=457500; Unidad="$/bins (420 kilos)"; Origen="Región de O'Higgins"; PKg=1089; KgUnidad=420 },
    @{ Row=69; Fecha=44539; Variedad="Artic Star"; Calidad="Segunda";  Volumen=16;  PMin=425000; PMax=430000; PProm=427500; Unidad="$/bins (420 kilos)"; Origen="Región de O'Higgins"; PKg=1018; KgUnidad=420 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value  = $mercadoId
    $ws.Cells.Item($row, 2).Value  = $mercado
    $ws.Cells.Item($row, 3).Value  = $region
    $ws.Cells.Item($row, 4).Value  = $r.Fecha
    $ws.Cells.Item($row, 5).Value  = $codreg
    $ws.Cells.Item($row, 6).Value  = $tipo
    $ws.Cells.Item($row, 7).Value  = $productoId
    $ws.Cells.Item($row, 8).Value  = $producto
    $ws.Cells.Item($row, 9).Value  = $categoriaId
    $ws.Cells.Item($row, 10).Value = $categoria
    $ws.Cells.Item($row, 11).Value = $r.Variedad
    $ws.Cells.Item($row, 12).Value = $r.Calidad
    $ws.Cells.Item($row, 13).Value = $r.Volumen
    $ws.Cells.Item($row, 14).Value = $r.PMin
    $ws.Cells.Item($row, 15).Value = $r.PMax
    $ws.Cells.Item($row, 16).Value = $r.PProm
    $ws.Cells.Item($row, 17).Value = $r.Unidad
    $ws.Cells.Item($row, 18).Value = $r.Origen
    $ws.Cells.Item($row, 19).Value = $r.PKg
    $ws.Cells.Item($row, 20).Value = $r.KgUnidad
}
